$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the 4 runs of the "Brief summary of conclusions/results..."
# paragraph into a single run. Re-setting identical text is a no-op in this
# engine, so we clear the range first and then set the final text so Word is
# forced to re-create a single run.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$r6.MoveEnd(1, -1)
$r6.Text = ""
$p6b = $d.Paragraphs(6)
$r6b = $p6b.Range
$r6b.MoveEnd(1, -1)
$r6b.Text = "Brief summary of conclusions/results drawn from experimental results"

# ---------------------------------------------------------------------------
# Change 2: highlight (yellow) the "Not much research..." sentence but leave
# the trailing period un-highlighted (this naturally splits it into its own
# run).
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$main12 = $p12.Range
$main12.MoveEnd(1, -2)
$main12.Font.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# Change 3: highlight (yellow) several whole paragraphs (including their
# paragraph mark, which records the highlight on the pPr/rPr too), skipping
# the "Provides a way to visualize..." paragraph.
# ---------------------------------------------------------------------------
$highlightParas = @(13, 14, 15, 17)
foreach ($idx in $highlightParas) {
    $p = $d.Paragraphs($idx)
    $full = $p.Range
    $full.Font.HighlightColorIndex = 7
}

# ---------------------------------------------------------------------------
# Change 4: move the "_GoBack" bookmark from the end of the "Compares
# results..." paragraph to the end of the "Experimental Design" paragraph.
#
# Adding a collapsed bookmark exactly at a paragraph's final text position
# has a quirk where it gets attached before the run instead of after it, so
# we work around it by temporarily appending a marker character, inserting
# the (now not-paragraph-final) bookmark, and then removing the marker.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$p19 = $d.Paragraphs(19)
$full19 = $p19.Range
$full19.MoveEnd(1, -1)
$endPos = $full19.End
$full19.InsertAfter("Z")

$p19b = $d.Paragraphs(19)
$bmRange = $p19b.Range
$bmRange.MoveEnd(1, -1)
$bmRange.SetRange($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$p19c = $d.Paragraphs(19)
$delRange = $p19c.Range
$delRange.MoveEnd(1, -1)
$delRange.SetRange($endPos, $endPos + 1)
$delRange.Text = ""
